# Append a freshly scraped job listing (2025-11-15 01:46 JST run) to the
# "ランサーズ" sheet. The scraper re-writes the whole sheet on every run, so
# every existing row's "取得日時" (retrieved-at) timestamp is refreshed, and
# the new listing is inserted in its sorted position at row 13 (pushing the
# old rows 13-18 down to 14-19).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newTimestamp = "2025-11-15 01:46:14"

# Make room for the new row; Excel shifts rows 13:18 down to 14:19 and
# carries their formatting (incl. hyperlinks) along with them.
$ws.Rows.Item(13).Insert()

# Fill in the new job listing that now occupies row 13.
$ws.Range("B13").Value = "【急募】TradingViewインジシグナルを用いたXAUUSD自動売買EA制作"
$ws.Range("C13").Value = "システム開発"
$ws.Range("D13").Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Range("E13").Value = "期限情報なし"
$ws.Range("F13").Value = "https://www.lancers.jp/work/detail/5434524"
$ws.Hyperlinks.Add($ws.Range("F13"), "https://www.lancers.jp/work/detail/5434524")
$ws.Range("F13").Style = "Hyperlink"
$ws.Range("G13").Value = 18

# Refresh the "取得日時" timestamp on every data row (the new row 13 as well
# as every row that existed before, rows 2:19 after the insert).
$ws.Range("A2:A19").Value = $newTimestamp
